# Update Sanju Samson / Rajasthan Royals innings stats (runs/balls/fours/sixes)
# on the single worksheet. Values are kept as text (matching the workbook's
# original "number stored as text" convention) rather than being converted
# to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, taken from the target diff.
$updates = @{
    "C2" = "4"
    "E2" = "1"

    "C4" = "85"
    "D4" = "42"
    "E4" = "4"
    "F4" = "7"

    "C5" = "26"
    "D5" = "25"
    "E5" = "3"

    "C6" = "5"
    "D6" = "9"
    "E6" = "0"
    "F6" = "0"

    "C7" = "0"
    "D7" = "3"
    "E7" = "0"

    "C8" = "8"
    "E8" = "1"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so the numeric-looking string isn't reinterpreted
    # as a number (keeps parity with the original t="str" cells).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
